$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
  "D2" = "27.487.77"
  "E2" = "  -1.49%  "
  "D3" = "1.574.63"
  "E3" = "  -3.61%  "
  "E4" = "  +0.43%  "
  "D5" = "205.27"
  "E5" = "  -2.96%  "
  "D6" = "0.501"
  "E6" = "  -3.49%  "
  "E7" = "  +0.48%  "
  "D8" = "21.98"
  "E8" = "  -6.09%  "
  "D9" = "0.250"
  "E9" = "  -2.83%  "
  "D10" = "0.0586"
  "E10" = "  -4.16%  "
  "D11" = "0.0863"
  "E11" = "  -2.17%  "
  "D12" = "1.797.20"
  "E12" = "  -3.63%  "
  "D13" = "1.581.98"
  "E13" = "  -3.20%  "
  "E14" = "  -5.03%  "
  "D15" = "0.527"
  "E15" = "  -6.57%  "
  "D16" = "27.468.85"
  "E16" = "  -1.61%  "
  "D17" = "62.52"
  "E17" = "  -4.47%  "
  "D18" = "215.77"
  "E18" = "  -5.92%  "
  "B19" = "ShibaInu"
  "C19" = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
  "D19" = "0.0₃0687"
  "E19" = "  -4.50%  "
  "B20" = "Chainlink"
  "C20" = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
  "D20" = "7.27"
  "E20" = "  -5.10%  "
  "D21" = "1.00"
  "E21" = "  +0.57%  "
  "E22" = "  -5.11%  "
  "E23" = "  -6.08%  "
  "D24" = "1.98"
  "D25" = "153.44"
  "E25" = "  -0.81%  "
  "E27" = "  -3.15%  "
  "E28" = "  -3.67%  "
  "E29" = "  -4.89%  "
  "E30" = "  -3.26%  "
  "D31" = "0.0461"
  "E31" = "  -4.11%  "
  "E32" = "  -5.95%  "
  "D33" = "1.357.24"
  "E33" = "  -2.73%  "
  "E34" = "  -6.02%  "
  "E35" = "  -5.67%  "
  "D36" = "0.967"
  "E36" = "  -5.22%  "
  "E37" = "  -1.21%  "
  "D38" = "0.0163"
  "E38" = "  -4.44%  "
  "D39" = "0.535"
  "E39" = "  -4.30%  "
  "D40" = "0.805"
  "E40" = "  -5.28%  "
  "D41" = "1.00"
  "E41" = "  +0.47%  "
  "D42" = "0.973"
  "E42" = "  -4.04%  "
  "D43" = "2.17"
  "E43" = "  +1.12%  "
  "E44" = "  -3.78%  "
  "D45" = "63.11"
  "E45" = "  -4.13%  "
  "D46" = "5.23"
  "E46" = "  -3.78%  "
  "D47" = "1.707.96"
  "E47" = "  -3.74%  "
  "D48" = "86.88"
  "E48" = "  -2.04%  "
  "B49" = "BabyDogeCoin"
  "C49" = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
  "D49" = "0.0₆0100"
  "E49" = "  -2.96%  "
  "B50" = "Algorand"
  "C50" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
  "D50" = "0.0962"
  "E50" = "  -5.37%  "
  "B51" = "Cronos"
  "C51" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
  "D51" = "0.0496"
  "E51" = "  -1.73%  "
}

foreach ($ref in $updates.Keys) {
  $ws.Range($ref).Value = "'" + $updates[$ref]
  $ws.Range($ref).Style = "Normal"
}
